# Template score sheet: add "Lần thứ" row and swap in new subject/column
# values (as opposed to the old "Học kì" / "Toán" / "Cuối kỳ" sample data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (old rows 6..16 shift down to 7..17,
# preserving the blank row that originally separated row 6 from "Danh sách").
$ws.Rows("6").Insert()

# Set the new/changed string values. The order below intentionally matches
# the order new shared strings get appended in the saved workbook so the
# resulting sharedStrings.xml layout matches the canonical output.
$ws.Range("B5").Value = "Miệng"
$ws.Range("A6").Value = "Lần thứ"
$ws.Range("B4").Value = "Ngữ Văn"
$ws.Range("A3").Value = "Học kỳ"

$ws.Range("A4").Value = "Môn học"
$ws.Range("A5").Value = "Cột điểm"

# Numeric value for the newly inserted row
$ws.Range("B6").Value = 1

# Match the saved selection state of the target workbook
$ws.Range("A3").Select()
